$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.494.03"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.17%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.261.78"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +6.65%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.80"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +4.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.74"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +7.71%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.253.21"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +6.70%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +5.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.06"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +9.21%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +5.95%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +5.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.11"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.31%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +5.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.778.95"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +6.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "558.47"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +12.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.583.85"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.260.29"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +6.69%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.18"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +6.04%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +5.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.751"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +8.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.89"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +9.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.69"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +7.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.17"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.44%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.30"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +18.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.01"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +8.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.28"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +6.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "27.92"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +6.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.79"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +6.33%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.18"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "570.77"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +8.28%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.45"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +7.01%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0456"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +11.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0871"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +7.81%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +7.55%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +11.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.219.75"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +10.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.71"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.283"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +14.34%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +10.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "26.71"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +6.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₃0562"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.06%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "126.25"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.38%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +4.04%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +8.08%  "
